$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 25
$ws.Cells.Item(4, 6).Value = 5854
$ws.Cells.Item(5, 6).Value = 67
$ws.Cells.Item(6, 6).Value = 2887
$ws.Cells.Item(7, 6).Value = 1253
$ws.Cells.Item(9, 6).Value = 383
$ws.Cells.Item(10, 6).Value = 427
$ws.Cells.Item(14, 6).Value = 154
$ws.Cells.Item(15, 6).Value = 4153
$ws.Cells.Item(16, 6).Value = 4153
$ws.Cells.Item(17, 6).Value = 89
$ws.Cells.Item(18, 6).Value = 79
$ws.Cells.Item(19, 6).Value = 88
$ws.Cells.Item(21, 6).Value = 188
$ws.Cells.Item(22, 6).Value = 51
$ws.Cells.Item(23, 6).Value = 6203
$ws.Cells.Item(24, 6).Value = 6203
$ws.Cells.Item(28, 6).Value = 423
$ws.Cells.Item(29, 6).Value = 195
$ws.Cells.Item(31, 6).Value = 5308
$ws.Cells.Item(32, 6).Value = 1591
$ws.Cells.Item(35, 6).Value = 5810
$ws.Cells.Item(38, 6).Value = 77
$ws.Cells.Item(39, 6).Value = 65
$ws.Cells.Item(40, 6).Value = 3877
$ws.Cells.Item(41, 6).Value = 88
$ws.Cells.Item(42, 6).Value = 66
$ws.Cells.Item(43, 6).Value = 12
$ws.Cells.Item(44, 6).Value = 2381
$ws.Cells.Item(49, 6).Value = 258
$ws.Cells.Item(50, 6).Value = 1347
$ws.Cells.Item(51, 6).Value = 14

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 182
$ws.Cells.Item(4, 6).Value = 23

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 1398

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1398
$ws.Cells.Item(3, 6).Value = 25
$ws.Cells.Item(4, 6).Value = 5854
$ws.Cells.Item(5, 6).Value = 67
$ws.Cells.Item(6, 6).Value = 2887
$ws.Cells.Item(7, 6).Value = 1253
$ws.Cells.Item(8, 6).Value = 427
$ws.Cells.Item(11, 6).Value = 182
$ws.Cells.Item(13, 6).Value = 154
$ws.Cells.Item(14, 6).Value = 4153
$ws.Cells.Item(15, 6).Value = 4153
$ws.Cells.Item(16, 6).Value = 89
$ws.Cells.Item(17, 6).Value = 79
$ws.Cells.Item(18, 6).Value = 88
$ws.Cells.Item(20, 6).Value = 188
$ws.Cells.Item(21, 6).Value = 51
$ws.Cells.Item(22, 6).Value = 6203
$ws.Cells.Item(23, 6).Value = 6203
$ws.Cells.Item(26, 6).Value = 423
$ws.Cells.Item(27, 6).Value = 195
$ws.Cells.Item(30, 6).Value = 5308
$ws.Cells.Item(31, 6).Value = 1591
$ws.Cells.Item(36, 6).Value = 5810
$ws.Cells.Item(39, 6).Value = 3877
$ws.Cells.Item(40, 6).Value = 66
$ws.Cells.Item(41, 6).Value = 12
$ws.Cells.Item(44, 6).Value = 2381
$ws.Cells.Item(49, 6).Value = 258
